$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.311.41"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "'1.874.01"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'0.7099"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "'241.96"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.07795"
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("D9").Value = "'0.3106"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'25.12"
$ws.Range("E10").Value = "  +1.68%  "
$ws.Range("D11").Value = "'0.08422"
$ws.Range("D12").Value = "'1.866.48"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "'5.238"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "'0.7162"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").Value = "'91.12"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "'29.320.35"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "'0.000008286"
$ws.Range("E17").Value = "  +2.48%  "
$ws.Range("D18").Value = "'6.078"
$ws.Range("E18").Value = "  +2.50%  "
$ws.Range("D19").Value = "'240.46"
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").Value = "'2.113.09"
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'7.743"
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("D25").Value = "'0.1596"
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'9.036"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'162.33"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("D28").Value = "'18.49"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "'1.298"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("D32").Value = "'4.324"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").Value = "'0.05376"
$ws.Range("E33").Value = "  +4.13%  "
$ws.Range("D34").Value = "'1.952"
$ws.Range("E34").Value = "  +2.08%  "
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("D36").Value = "'0.7503"
$ws.Range("E36").Value = "  -2.94%  "
$ws.Range("D37").Value = "'2.688"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "'0.01886"
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("D39").Value = "'1.235.18"
$ws.Range("E39").Value = "  +6.44%  "
$ws.Range("D40").Value = "'2.729"
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("D42").Value = "'0.8944"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").Value = "'109.01"
$ws.Range("E43").Value = "  +5.76%  "
$ws.Range("D44").Value = "'72.41"
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "'2.013.72"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("E47").Value = "  +6.23%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.5200"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.796"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").Value = "'9.450"
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("D51").Value = "'0.4348"
$ws.Range("E51").Value = "  +1.23%  "
